$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra data rows (3-8) - only the header row and the first
# data row (BillTo / Coordinates PDF row) remain.
$ws.Rows("3:8").Delete()

# The "PASS" status for the remaining row is cleared out.
$ws.Range("C2").ClearContents()

# Narrow columns A and B slightly.
$ws.Columns("A").ColumnWidth = 64.65
$ws.Columns("B").ColumnWidth = 50.5

# Selection moves to A2.
$ws.Range("A2").Select()
